# Scheduled-runner update: refresh cached market-board price snapshots
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ),
# the derived Leve price totals (K/L) and profit figures (M/N) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 10871217
$ws.Range("I15").Value = 10871217
$ws.Range("K15").Value = 32613651
$ws.Range("M15").Value = -32613482

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 240925.1
$ws.Range("J17").Value = 240925.1
$ws.Range("L17").Value = 722775.3
$ws.Range("N17").Value = -723111.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1005.1724
$ws.Range("I98").Value = 1110.8695
$ws.Range("J98").Value = 600
$ws.Range("K98").Value = 1110.8695
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = 387.1305
$ws.Range("N98").Value = -3596

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1005.1724
$ws.Range("I122").Value = 1110.8695
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3332.6085
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -882.6085000000003
$ws.Range("N122").Value = -6700

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 314519.38
$ws.Range("I132").Value = 2101.1667
$ws.Range("K132").Value = 6303.500100000001
$ws.Range("M132").Value = -3773.500100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2435.84
$ws.Range("J138").Value = 2793
$ws.Range("L138").Value = 8379
$ws.Range("N138").Value = -18659

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 49952.855
$ws.Range("J139").Value = 49952.855
$ws.Range("L139").Value = 49952.855
$ws.Range("N139").Value = -60232.855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1928.6
$ws.Range("I45").Value = 1660.75
$ws.Range("K45").Value = 1660.75
$ws.Range("M45").Value = -1283.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 900.6957
$ws.Range("I74").Value = 855.9487
$ws.Range("K74").Value = 855.9487
$ws.Range("M74").Value = 18.05129999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 900.6957
$ws.Range("I77").Value = 855.9487
$ws.Range("K77").Value = 4279.7435
$ws.Range("M77").Value = 88.25649999999951

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 39333
$ws.Range("J92").Value = 39333
$ws.Range("L92").Value = 39333
$ws.Range("N92").Value = -44325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 17850
$ws.Range("J121").Value = 17850
$ws.Range("L121").Value = 17850
$ws.Range("N121").Value = -21344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1555.75
$ws.Range("I122").Value = 1555.75
$ws.Range("K122").Value = 4667.25
$ws.Range("M122").Value = -2217.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1787.72
$ws.Range("I132").Value = 1204.5333
$ws.Range("K132").Value = 3613.5999
$ws.Range("M132").Value = -1083.5999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 12778
$ws.Range("I102").Value = 12778
$ws.Range("K102").Value = 12778
$ws.Range("M102").Value = -9533

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6014.2144
$ws.Range("I105").Value = 5019.9
$ws.Range("J105").Value = 8500
$ws.Range("K105").Value = 5019.9
$ws.Range("L105").Value = 8500
$ws.Range("M105").Value = -3272.9
$ws.Range("N105").Value = -11994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28541.8
$ws.Range("I134").Value = 3885.1365
$ws.Range("J134").Value = 58677.723
$ws.Range("K134").Value = 11655.4095
$ws.Range("L134").Value = 176033.169
$ws.Range("M134").Value = -9120.4095
$ws.Range("N134").Value = -181103.169

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 49780
$ws.Range("J138").Value = 49780
$ws.Range("L138").Value = 49780
$ws.Range("N138").Value = -60060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3726.4
$ws.Range("I31").Value = 3726.4
$ws.Range("K31").Value = 3726.4
$ws.Range("M31").Value = -3431.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3726.4
$ws.Range("I34").Value = 3726.4
$ws.Range("K34").Value = 3726.4
$ws.Range("M34").Value = -3524.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2717.3333
$ws.Range("I132").Value = 2021.9412
$ws.Range("J132").Value = 3899.5
$ws.Range("K132").Value = 6065.8236
$ws.Range("L132").Value = 11698.5
$ws.Range("M132").Value = -3535.8236
$ws.Range("N132").Value = -16758.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 53167.273
$ws.Range("J140").Value = 53167.273
$ws.Range("L140").Value = 53167.273
$ws.Range("N140").Value = -63527.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1277.1177
$ws.Range("I5").Value = 963.29034
$ws.Range("K5").Value = 2889.87102
$ws.Range("M5").Value = -2777.87102

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1277.1177
$ws.Range("I135").Value = 963.29034
$ws.Range("K135").Value = 8669.61306
$ws.Range("M135").Value = -6134.61306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4669.222
$ws.Range("I70").Value = 4100
$ws.Range("J70").Value = 5807.6665
$ws.Range("K70").Value = 4100
$ws.Range("L70").Value = 5807.6665
$ws.Range("M70").Value = -3830
$ws.Range("N70").Value = -6347.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4669.222
$ws.Range("I73").Value = 4100
$ws.Range("J73").Value = 5807.6665
$ws.Range("K73").Value = 4100
$ws.Range("L73").Value = 5807.6665
$ws.Range("M73").Value = -3164
$ws.Range("N73").Value = -7679.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2434.375
$ws.Range("I122").Value = 1794.85
$ws.Range("K122").Value = 5384.549999999999
$ws.Range("M122").Value = -2934.549999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 22329.428
$ws.Range("J136").Value = 22329.428
$ws.Range("L136").Value = 66988.284
$ws.Range("N136").Value = -72088.284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 35750
$ws.Range("J138").Value = 35750
$ws.Range("L138").Value = 35750
$ws.Range("N138").Value = -46030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 19000
$ws.Range("J98").Value = 19000
$ws.Range("L98").Value = 19000
$ws.Range("N98").Value = -24990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 10000000
$ws.Range("J116").Value = 10000000
$ws.Range("L116").Value = 10000000
$ws.Range("N116").Value = -10009178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5183.6763
$ws.Range("I122").Value = 7563.1113
$ws.Range("K122").Value = 22689.3339
$ws.Range("M122").Value = -20239.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3098.3447
$ws.Range("I132").Value = 2399.625
$ws.Range("J132").Value = 3958.3076
$ws.Range("K132").Value = 7198.875
$ws.Range("L132").Value = 11874.9228
$ws.Range("M132").Value = -4668.875
$ws.Range("N132").Value = -16934.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2715.2666
$ws.Range("I136").Value = 1324.0526
$ws.Range("J136").Value = 5118.273
$ws.Range("K136").Value = 3972.1578
$ws.Range("L136").Value = 15354.819
$ws.Range("M136").Value = -1422.1578
$ws.Range("N136").Value = -20454.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 32250
$ws.Range("J137").Value = 32250
$ws.Range("L137").Value = 32250
$ws.Range("N137").Value = -42450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 36966.668
$ws.Range("J138").Value = 36966.668
$ws.Range("L138").Value = 36966.668
$ws.Range("N138").Value = -47246.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 2602
$ws.Range("J101").Value = 2602
$ws.Range("L101").Value = 2602
$ws.Range("N101").Value = -9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1293.5625
$ws.Range("I122").Value = 1250.5
$ws.Range("K122").Value = 3751.5
$ws.Range("M122").Value = -1301.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1667.6207
$ws.Range("I136").Value = 1602.1786
$ws.Range("K136").Value = 4806.5358
$ws.Range("M136").Value = -2256.5358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 45600
$ws.Range("J138").Value = 45600
$ws.Range("L138").Value = 45600
$ws.Range("N138").Value = -55880
